$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the columns that are no longer part of the simplified "penduduk" import
# template. Work from right to left so earlier deletions don't shift the
# column letters we still need to remove.
#
# Original header layout (A1:P1):
#   A Nama
#   B NIK
#   C Jenis Kelamin
#   D Tempat Lahir            <- remove
#   E Tanggal Lahir
#   F Agama                   <- remove
#   G Status Pendidikan Terakhir
#   H Pekerjaan
#   I Golongan Darah          <- remove
#   J Status Perkawinan       <- remove
#   K Tanggal Perkawinan      <- remove
#   L Kewarganegaraan         <- remove
#   M Nomor Paspor            <- remove
#   N Nomor Kitap             <- remove
#   O Alamat                  <- remove
#   P Desa
#
# Resulting header layout (A1:G1):
#   A Nama
#   B NIK
#   C Jenis Kelamin
#   D Tanggal Lahir
#   E Status Pendidikan Terakhir
#   F Pekerjaan
#   G Desa

$ws.Range("O1").EntireColumn.Delete()
$ws.Range("N1").EntireColumn.Delete()
$ws.Range("M1").EntireColumn.Delete()
$ws.Range("L1").EntireColumn.Delete()
$ws.Range("K1").EntireColumn.Delete()
$ws.Range("J1").EntireColumn.Delete()
$ws.Range("I1").EntireColumn.Delete()
$ws.Range("F1").EntireColumn.Delete()
$ws.Range("D1").EntireColumn.Delete()

# Set the column widths that Excel records once the sheet was tidied up.
$ws.Range("B1").EntireColumn.ColumnWidth = 20.28515625
$ws.Range("C1").EntireColumn.ColumnWidth = 24.42578125
$ws.Range("D1").EntireColumn.ColumnWidth = 26.7109375
$ws.Range("E1").EntireColumn.ColumnWidth = 27
$ws.Range("F1").EntireColumn.ColumnWidth = 21.42578125

# Move the active selection off the header row, matching the saved view state.
$ws.Range("E15").Select()
